$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107 (shifts old rows 107+ down by one, formulas referencing
# those rows get auto-adjusted by Excel, matching the B109->B110, B110->B111 shift).
$ws.Rows(107).Insert()

# Row 106 is no longer the "latest" (highlighted) row, so it should pick up the
# plain/"Neutral" formatting used by the rows above it (copy formats only).
$ws.Range("A105:I105").Copy()
$ws.Range("A106:I106").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The "daily rate to achieve June 20 target" note moves from I106 to the new I107.
$ws.Cells.Item(106, 9).ClearContents()

# New data row 107 - latest day of data (keeps the "Good"/highlighted style that
# Excel's row-insert already carried down from the old row 106).
$ws.Cells.Item(107, 1).Value = 44292
$ws.Cells.Item(107, 2).Value = 3299
$ws.Cells.Item(107, 3).Formula = "=(AVERAGE(B101:B107))"
$ws.Cells.Item(107, 4).Formula = "=(D106-B107)"
$ws.Cells.Item(107, 5).Formula = "=E106+B107"
$ws.Cells.Item(107, 6).Formula = "=D107/C107"
$ws.Cells.Item(107, 7).Formula = "=A107+F107"
$ws.Cells.Item(107, 8).Formula = "=D107/84"
$ws.Cells.Item(107, 9).Value = "daily rate to achieve June 20 target"

# Restore the sheet view/selection to what Excel would leave it at after this edit.
$ws.Application.Goto($ws.Range("G112"), $true)
$ws.Range("G112").Select()
